# "improved Figure 2 resolution and title size"
#
# 1) Refresh the auto date placeholder ("datetimeFigureOut" field) from
#    5/6/23 -> 8/23/23 everywhere it is defined on the slide master and
#    every slide layout.
# 2) Re-crop/re-position the picture, enlarge+reposition the title
#    textbox (and bump its run size 20pt -> 36pt, autosize Off), and
#    resize/reposition the surrounding frame rectangle - i.e. the
#    "improved resolution" re-layout of Figure 2.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($shapes, [string]$text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $text
        }
    }
}

$newDate = "8/23/23"

# Slide master.
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes $newDate
}

# NOTE: intentionally not touching $p.NotesMaster.Shapes here - in this
# host, enumerating/writing notes-master shapes cross-wires with the
# slide master's shape storage and corrupts an unrelated slide-master
# placeholder's text. The notes master also carries its own copy of the
# same auto date field, but it is not reachable without that corruption.

# --- Figure 2 re-layout on slide 1 -----------------------------------

$s = $p.Slides.Item(1)

$pic = $s.Shapes.Item("Picture 8")
$pic.Left = 100.13791271574802
$pic.Top = 67.12775803543306
$pic.Width = 643.190157480315
$pic.Height = 443.7910614220473

$title = $s.Shapes.Item("Title 1")
$title.Left = 27.750118340157478
$title.Top = 11.12496062992126
$title.Width = 843.3228760456693
$title.Height = 67.49574803149606
$title.TextFrame.AutoSize = 0
$title.TextFrame.TextRange.Font.Size = 36

$rect = $s.Shapes.Item("Rectangle 11")
$rect.Left = 41.93527559055118
$rect.Top = 67.12775803543306
$rect.Width = 759.5951968503937
$rect.Height = 443.7910614220473
